$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Mise" (stake) base value used throughout the sheet
$ws.Range("E3").Value = 0.0042

# Correct the last row counter (C11) from 8 to 9
$ws.Range("C11").Value = 9

# Correct the target dollar amount used to compute time required
$ws.Range("I7").Formula = "=(10/I5)/60"

# Update the active selection left by the author when saving
$ws.Range("H15").Select()
